# Update simple accessioning forms
# Insert a new "Tags:" column between the existing "Sequencing Date:" (O)
# and "Files:" (old P, now shifted to Q) columns on the single worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Inserting at column P shifts P..T -> Q..U, carrying over the formatting
# (style/width) of the column immediately to the left (O), matching how
# Excel performs a manual "Insert Column" from the right-click menu.
$ws.Columns("P").Insert()

# Populate the new header cell with the new shared string.
$ws.Range("P1").Value = "Tags:"

# Match the width of the newly inserted column to its left-hand neighbor
# (O) so the column isn't left at the sheet default width.
$ws.Columns("P").ColumnWidth = $ws.Columns("O").ColumnWidth

# Reflect the updated selection/viewport left in the saved workbook.
$ws.Range("O4").Select() | Out-Null
